$d = $word.ActiveDocument

# Replace the full text of every paragraph whose text equals $oldText with
# $newText, while preserving the paragraph's existing structure (pPr, any
# empty/placeholder runs, run formatting, etc.) exactly as-is. A plain
# Find.Execute / Range.Text replace would work for the visible text, but this
# engine (like Word itself) re-normalizes/collapses a paragraph's runs
# whenever its text is edited in place, which silently drops zero-length
# sibling runs. Rebuilding the paragraph from its own WordOpenXML and
# re-inserting it keeps that paragraph-internal markup untouched.
function Replace-ParagraphText($oldText, $newText) {
    $targets = @()
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $oldText) {
            $targets += $p
        }
    }
    if ($targets.Count -eq 0) {
        throw "Replace-ParagraphText: no paragraph found with text: $oldText"
    }
    foreach ($target in $targets) {
        $xml = $target.Range.WordOpenXML
        if ($xml -notmatch '<w:body>(.*?)</w:body>') {
            throw "Replace-ParagraphText: could not locate w:body in WordOpenXML"
        }
        $body = $matches[1]
        if ($body -match '^(<w:p\b[^>]*>.*?</w:p>|<w:p\b[^>]*/>)') {
            $pxml = $matches[1]
        } else {
            throw "Replace-ParagraphText: could not isolate paragraph xml"
        }
        # Drop the fresh w14:paraId/rsid* identifiers the WordOpenXML export
        # mints on round-trip so InsertXML doesn't stamp them onto the
        # paragraph mark.
        $pxml = $pxml -replace ' w14:paraId="[0-9A-Fa-f]*"', ''
        $pxml = $pxml -replace ' w14:textId="[0-9A-Fa-f]*"', ''
        $pxml = $pxml -replace ' w:rsidR="[0-9A-Fa-f]*"', ''
        $pxml = $pxml -replace ' w:rsidRDefault="[0-9A-Fa-f]*"', ''
        $pxml = $pxml -replace ' w:rsidRPr="[0-9A-Fa-f]*"', ''
        $pxml = $pxml -replace ' w:rsidP="[0-9A-Fa-f]*"', ''

        $escapedOld = [regex]::Escape($oldText)
        $newXml = $pxml -replace (">" + $escapedOld + "</w:t>"), (">" + $newText + "</w:t>")
        if ($newXml -eq $pxml) {
            throw "Replace-ParagraphText: text substitution did not match inside: $pxml"
        }
        $target.Range.InsertXML($newXml)
    }
}

# Title / heading text (appears twice: Heading1 and the bold run near the end)
Replace-ParagraphText "Play Book of Vikings Free - Review of Exciting Norse Mythology-Inspired Slot" "Play Book of Vikings for Free Online"

# "What we like" bullet list
Replace-ParagraphText "Highly immersive Norse mythology theme" "Classic 5x3 grid with 10 fixed paylines"
Replace-ParagraphText "Autoplay feature with win/loss limits for controlled gameplay" "Autoplay function with win/loss limits"
Replace-ParagraphText "Substantial payouts with up to 200x total bet" "Bonus Round with 10 free spins and expanding special symbol"
Replace-ParagraphText "Suitable for players of all levels with intuitive design" "Visually engaging design inspired by Norse mythology"

# "What we don't like" bullet list
Replace-ParagraphText "High volatility may not appeal to all players" "High volatility may not appeal to players who prefer frequent wins"
Replace-ParagraphText "Only 10 fixed paylines" "No progressive jackpot feature"

# Closing italic summary paragraph
Replace-ParagraphText "Read our review of Book of Vikings, an exciting online slot game inspired by Norse mythology. Play for free and potentially win big with high volatility gameplay." "Read our review of Book of Vikings and play for free online. Enjoy immersive gameplay and high volatility."
